# Automatische test-sync: 2025-06-22 19:07:50
#
# Adds a new incoming-mail log entry (row 37) to the "Logs" sheet and
# refreshes the category-count table on the "Dashboard" sheet to reflect it.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Logs sheet: append the new row at the bottom of the table.
# ---------------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$newRow = 37
$logs.Cells.Item($newRow, 1).Value = "Sollicitatie salesfunctie"
$logs.Cells.Item($newRow, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($newRow, 3).Value = "Hierbij mijn sollicitatie voor de salesfunctie. CV in bijlage."
$logs.Cells.Item($newRow, 4).Value = "Sollicitatie / Vacature"
$logs.Cells.Item($newRow, 5).Value = "Beste sollicitant,`nDank voor het sturen van je sollicitatie voor de salesfunctie. We zullen je CV zorgvuldig bekijken en zo spoedig mogelijk contact met je opnemen over de vervolgstappen. Mocht je in de tussentijd vragen hebben, dan horen we het graag.`nMet vriendelijke groet,`n[Naam bedrijf]"
$logs.Cells.Item($newRow, 6).Value = "2025-06-22 19:07:11"
$logs.Cells.Item($newRow, 7).Value = "Ja"

# Re-fit the new row's height back to the sheet default — entering the
# multi-line reply text would otherwise leave an explicit custom row
# height behind, which the source row never had.
$logs.Rows.Item($newRow).EntireRow.AutoFit()

# Extend the two conditional-formatting rule ranges (Categorie / Beantwoord)
# so they keep covering the whole table, now through row 37.
$catRules = $logs.Range("D2:D36").FormatConditions
for ($i = 1; $i -le $catRules.Count; $i++) {
    $catRules.Item($i).ModifyAppliesToRange($logs.Range("D2:D37"))
}

$answeredRules = $logs.Range("G2:G36").FormatConditions
for ($i = 1; $i -le $answeredRules.Count; $i++) {
    $answeredRules.Item($i).ModifyAppliesToRange($logs.Range("G2:G37"))
}

# ---------------------------------------------------------------------------
# 2. Dashboard sheet: update the category summary table.
#
#    The new "Sollicitatie / Vacature" entry raises that category's count
#    from 2 to 3, which re-sorts the (count desc) summary table: the
#    "Afmelding / Nieuwsbrief" / "Offerte / Prijsaanvraag" / "Overig" /
#    "Juridisch / Contract" / "Sollicitatie / Vacature" /
#    "Openingstijden / Locatie" block gets re-ordered accordingly.
# ---------------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Cells.Item(6, 1).Value = "Offerte / Prijsaanvraag"
$dash.Cells.Item(6, 2).Value = 3

$dash.Cells.Item(7, 1).Value = "Sollicitatie / Vacature"
$dash.Cells.Item(7, 2).Value = 3

$dash.Cells.Item(8, 1).Value = "Afmelding / Nieuwsbrief"
$dash.Cells.Item(8, 2).Value = 3

$dash.Cells.Item(9, 1).Value = "Overig"
$dash.Cells.Item(9, 2).Value = 2

$dash.Cells.Item(10, 1).Value = "Openingstijden / Locatie"
$dash.Cells.Item(10, 2).Value = 2

$dash.Cells.Item(11, 1).Value = "Juridisch / Contract"
$dash.Cells.Item(11, 2).Value = 2
